$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- 1. Remove the now-aged rows (old 9-16) entirely, shrinking the used range ---
$ws.Range("A9:H16").Delete() | Out-Null

# --- 2. Refresh rows 2-8 with the newly scraped listings ---
# Row 2
$ws.Range("A2").Value = '2025-12-19 06:29:19'
$ws.Range("B2").Value = 'EC×AIプロダクト/業務改善リード'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5450024'
$ws.Range("G2").Value = 338
$ws.Range("H2").Value = '🔥AI,Ai ◇業務改善'

# Row 3
$ws.Range("A3").Value = '2025-12-19 06:29:19'
$ws.Range("B3").Value = '製造業のR&D支援!「プロセスデータ解析」「音響異常検知」のAIエンジニア募集'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5439165'
$ws.Range("G3").Value = 303
$ws.Range("H3").Value = '🔥AI,Ai'

# Row 4
$ws.Range("A4").Value = '2025-12-19 06:29:19'
$ws.Range("B4").Value = '【急募】AWSスクレイピングツールの開発を依頼したいです!'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5457255'
$ws.Range("G4").Value = 165
$ws.Range("H4").Value = '◆ツール,開発'

# Row 5
$ws.Range("A5").Value = '2025-12-19 06:29:19'
$ws.Range("B5").Value = '【急募】飲食店予約サイトの制作と将来的なアプリ化(アプリ化の際は別契約)'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5457089'
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = '◇アプリ'

# Row 6
$ws.Range("A6").Value = '2025-12-19 06:29:19'
$ws.Range("B6").Value = '【外国人大歓迎】【急募】ECツールの保守・バグ修正・機能追加エンジニア募集'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5457026'
$ws.Range("G6").Value = 68
$ws.Range("H6").Value = '◆ツール'

# Row 7
$ws.Range("A7").Value = '2025-12-19 06:29:19'
$ws.Range("B7").Value = '【急募】PHPによる申請サイト構築支援!'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5457023'
$ws.Range("G7").Value = 58
$ws.Range("H7").Value = '◇サイト ○PHP'

# Row 8
$ws.Range("A8").Value = '2025-12-19 06:29:19'
$ws.Range("B8").Value = '【急募】Kintoneでの請求書自動発行システム構築依頼'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5457134'
$ws.Range("G8").Value = 28

# Row 8 has no "skill summary" note this time - make sure H8 is empty
$ws.Range("H8").ClearContents()

# --- 3. Rebuild the URL hyperlinks so only F2:F8 carry live links ---
$ws.Range("A1:H100").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5450024') | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5439165') | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5457255') | Out-Null
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5457089') | Out-Null
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5457026') | Out-Null
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5457023') | Out-Null
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5457134') | Out-Null
$ws.Range("F8").Style = "Hyperlink"

# --- 4. Narrow columns B and H to match the new content widths ---
# (ColumnWidth round-trips through this host with a constant +5/6 char
#  padding vs. the stored <col width>, so pre-subtract it to land exactly
#  on the target widths of 41 / 14.)
$ws.Columns.Item(2).ColumnWidth = 41 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 14 - 0.8333333333333334

